$wb = $excel.ActiveWorkbook

# --- Sheet "Variables": fix variable-name typos ---
$wsVars = $wb.Worksheets.Item("Variables")
$wsVars.Range("B7").Value = "m_berufab"
$wsVars.Range("B8").Value = "v_berufab"

# --- Sheet "Categories": update the "variable" column to match the renamed variables ---
$wsCats = $wb.Worksheets.Item("Categories")
for ($r = 14; $r -le 21; $r++) {
    $wsCats.Cells.Item($r, 1).Value = "m_berufab"
}
for ($r = 22; $r -le 29; $r++) {
    $wsCats.Cells.Item($r, 1).Value = "v_berufab"
}
